$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 134 (shifts existing rows 134-161 down to 135-162)
$ws.Rows.Item(134).EntireRow.Insert()

# Populate the newly inserted row 134 with the latest weekly data point
$ws.Range("A134").Value = 4
$ws.Range("B134").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C134").Value = "Los Lagos"
$ws.Range("D134").Value = 45135
$ws.Range("E134").Value = 10
$ws.Range("F134").Value = 100112031
$ws.Range("G134").Value = "Poroto verde"
$ws.Range("H134").Value = "Magnum"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 40
$ws.Range("K134").Value = 32000
$ws.Range("L134").Value = 32000
$ws.Range("M134").Value = 32000
$ws.Range("N134").Value = "$/malla 25 kilos"
$ws.Range("O134").Value = "Perú"
$ws.Range("P134").Value = 1280
$ws.Range("Q134").Value = 25
$ws.Range("R134").Value = "Hortaliza"
